# Apply updates to rows 16 and 17 (Indonesian PoS / PUD test results)
# following the change of Indonesian PoS and PUD tests data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 ---
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = 0.9984281898500103
$ws.Range("I16").Value = 0.9263399936568347
$ws.Range("J16").Value = "2:55:24"
$ws.Range("K16").Value = 4477
$ws.Range("L16").Value = 30.17176680813045
$ws.Range("M16").Value = 4477

# --- Row 17 ---
$ws.Range("E17").Value = 16
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 0.9844664773251581
$ws.Range("I17").Value = 0.9276086267047257
$ws.Range("J17").Value = "1:42:48"
$ws.Range("K17").Value = 4477
$ws.Range("L17").Value = 30.61000670091579
$ws.Range("M17").Value = 4477
